# The weekly Perejil price series gained a new observation.
# A new row is inserted at row 93 (pushing the existing rows 93-147 down to
# 94-148, each keeping its original data), and the new row 93 is populated
# with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 93:147 down to 94:148, inserting a blank row at 93.
$ws.Rows(93).Insert()

# Populate the newly inserted row 93 with this week's record.
$ws.Cells.Item(93, 1).Value  = 8
$ws.Cells.Item(93, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(93, 3).Value  = "Coquimbo"
$ws.Cells.Item(93, 4).Value  = 44719
$ws.Cells.Item(93, 5).Value  = 4
$ws.Cells.Item(93, 6).Value  = 100112044
$ws.Cells.Item(93, 7).Value  = "Perejil"
$ws.Cells.Item(93, 8).Value  = "Sin especificar"
$ws.Cells.Item(93, 9).Value  = "Primera"
$ws.Cells.Item(93, 10).Value = 3200
$ws.Cells.Item(93, 11).Value = 1500
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = 1750
$ws.Cells.Item(93, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(93, 15).Value = "Provincia del Elqu" + [char]0xED
$ws.Cells.Item(93, 16).Value = 1167
$ws.Cells.Item(93, 17).Value = 1.5
$ws.Cells.Item(93, 18).Value = "Hortaliza"
